$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (date advanced by one day)
$ws.Name = "Through 2022-03-30"

# Update the "March" row label to reflect the new date
$ws.Range("A4").Value = "March (through 03-30)"

# Update March row (row 4) values
$ws.Range("B4").Value = 29
$ws.Range("D4").Value = 57
$ws.Range("E4").Value = 60
$ws.Range("F4").Value = 30
$ws.Range("G4").Value = 56
$ws.Range("H4").Value = 79
$ws.Range("I4").Value = 130

# Update Total row (row 5) values
$ws.Range("B5").Value = 66
$ws.Range("D5").Value = 188
$ws.Range("E5").Value = 197
$ws.Range("F5").Value = 109
$ws.Range("G5").Value = 197
$ws.Range("H5").Value = 421
$ws.Range("I5").Value = 430
